$d = $word.ActiveDocument

# Locate the sentence that needs to be rewritten (spans the two existing
# runs "This document outline" + "s").
$searchRng = $d.Content
$found = $searchRng.Find.Execute("This document outlines")

if ($found) {
    $rng = $d.Range($searchRng.Start, $searchRng.End)
    $w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
    $xml = '<w:p xmlns:w="' + $w + '">' +
           '<w:r><w:t>Th</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">is </w:t></w:r>' +
           '<w:r><w:t>document</w:t></w:r>' +
           '<w:r><w:t>ation</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">covers the project overview, work-breakdown structure, activity definition and Gantt chart of the NSW Traffic Penalty software. </w:t></w:r>' +
           '</w:p>'
    $rng.InsertXML($xml)
}
